$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "25.511.59"
$ws.Range("E2").Value = "  -1.44%  "

# Row 3
$ws.Range("D3").Value = "1.571.05"
$ws.Range("E3").Value = "  -3.06%  "

# Row 4
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.60%  "

# Row 5
$ws.Range("D5").Value = "'207.43"
$ws.Range("E5").Value = "  -2.57%  "

# Row 6
$ws.Range("D6").Value = "'0.998"
$ws.Range("E6").Value = "  -0.70%  "

# Row 7
$ws.Range("D7").Value = "'0.480"
$ws.Range("E7").Value = "  -3.87%  "

# Row 8
$ws.Range("D8").Value = "'0.0616"
$ws.Range("E8").Value = "  -0.16%  "

# Row 9
$ws.Range("D9").Value = "'0.245"
$ws.Range("E9").Value = "  -1.96%  "

# Row 10
$ws.Range("D10").Value = "'18.00"
$ws.Range("E10").Value = "  -2.70%  "

# Row 11
$ws.Range("D11").Value = "'0.0782"
$ws.Range("E11").Value = "  -0.93%  "

# Row 12
$ws.Range("D12").Value = "1.781.17"
$ws.Range("E12").Value = "  -3.48%  "

# Row 13
$ws.Range("D13").Value = "1.561.19"
$ws.Range("E13").Value = "  -3.64%  "

# Row 14
$ws.Range("D14").Value = "'4.00"
$ws.Range("E14").Value = "  -3.44%  "

# Row 15
$ws.Range("D15").Value = "'0.510"
$ws.Range("E15").Value = "  -2.64%  "

# Row 16
$ws.Range("D16").Value = "25.401.82"
$ws.Range("E16").Value = "  -1.86%  "

# Row 17
$ws.Range("D17").Value = "0.0₃0719"
$ws.Range("E17").Value = "  -2.58%  "

# Row 18
$ws.Range("D18").Value = "'59.70"
$ws.Range("E18").Value = "  -3.06%  "

# Row 19
$ws.Range("E19").Value = "  -0.28%  "

# Row 20
$ws.Range("D20").Value = "'188.35"
$ws.Range("E20").Value = "  -1.84%  "

# Row 21
$ws.Range("D21").Value = "'4.16"
$ws.Range("E21").Value = "  -2.02%  "

# Row 22
$ws.Range("D22").Value = "'9.35"
$ws.Range("E22").Value = "  -1.61%  "

# Row 23
$ws.Range("D23").Value = "'5.90"
$ws.Range("E23").Value = "  -2.14%  "

# Row 24
$ws.Range("D24").Value = "'0.132"
$ws.Range("E24").Value = "  -2.69%  "

# Row 25
$ws.Range("D25").Value = "'141.00"
$ws.Range("E25").Value = "  -2.11%  "

# Row 26
$ws.Range("E26").Value = "  -0.47%  "

# Row 27
$ws.Range("E27").Value = "  -2.24%  "

# Row 28
$ws.Range("D28").Value = "'15.06"
$ws.Range("E28").Value = "  -0.87%  "

# Row 29
$ws.Range("D29").Value = "'6.42"
$ws.Range("E29").Value = "  -3.81%  "

# Row 30
$ws.Range("D30").Value = "'1.15"
$ws.Range("E30").Value = "  -6.38%  "

# Row 31
$ws.Range("D31").Value = "'0.0468"
$ws.Range("E31").Value = "  -1.86%  "

# Row 32
$ws.Range("D32").Value = "'3.08"
$ws.Range("E32").Value = "  -1.60%  "

# Row 33
$ws.Range("D33").Value = "'3.00"
$ws.Range("E33").Value = "  -3.13%  "

# Row 34
$ws.Range("E34").Value = "  +0.41%  "

# Row 35
$ws.Range("E35").Value = "  -4.53%  "

# Row 36
$ws.Range("D36").Value = "1.092.98"
$ws.Range("E36").Value = "  -2.97%  "

# Row 37
$ws.Range("B37").Value = "MXToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D37").Value = "'2.34"
$ws.Range("E37").Value = "  -1.17%  "

# Row 38
$ws.Range("B38").Value = "PaxDollar"
$ws.Range("C38").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D38").Value = "'0.998"
$ws.Range("E38").Value = "  -0.65%  "

# Row 39
$ws.Range("D39").Value = "'0.500"
$ws.Range("E39").Value = "  -2.05%  "

# Row 40
$ws.Range("E40").Value = "  -2.45%  "

# Row 41
$ws.Range("D41").Value = "'0.780"
$ws.Range("E41").Value = "  -7.07%  "

# Row 42
$ws.Range("D42").Value = "'0.801"
$ws.Range("E42").Value = "  +6.64%  "

# Row 43
$ws.Range("D43").Value = "'93.25"
$ws.Range("E43").Value = "  -4.99%  "

# Row 44
$ws.Range("D44").Value = "'5.13"
$ws.Range("E44").Value = "  +2.00%  "

# Row 45
$ws.Range("D45").Value = "1.697.25"
$ws.Range("E45").Value = "  -3.34%  "

# Row 46
$ws.Range("D46").Value = "0.0₆0105"
$ws.Range("E46").Value = "  -5.88%  "

# Row 47
$ws.Range("D47").Value = "'1.49"
$ws.Range("E47").Value = "  -1.49%  "

# Row 48
$ws.Range("D48").Value = "'52.90"
$ws.Range("E48").Value = "  -2.18%  "

# Row 49
$ws.Range("D49").Value = "'0.0507"
$ws.Range("E49").Value = "  -2.57%  "

# Row 50
$ws.Range("D50").Value = "'0.405"
$ws.Range("E50").Value = "  -1.71%  "

# Row 51
$ws.Range("D51").Value = "'0.998"
$ws.Range("E51").Value = "  -0.64%  "

